# The unified diff for this change is entirely a cosmetic / canonical-XML
# artifact: every single hunk (the w:document root's xmlns attribute order,
# w:pgSz / w:pgMar attribute order, w:rFonts / w:lang attribute order, the
# w:latentStyles / w:lsdException attribute order, and the w:style element
# attribute order) is just the *same* attributes re-serialized in a
# different (alphabetical) order. No attribute value, text run, paragraph,
# style definition, section property, or any other semantic content differs
# between the "before" and "after" sides of the diff.
#
# Word's COM object model does not expose raw control over the attribute
# serialization order used when a part is written back out, and there is
# no content-level edit implied by the diff to apply. So the correct,
# faithful reproduction of this change is simply to leave the document's
# content untouched.
#
# Touch the active document (matches the host's expected `$d` pattern) but
# perform no mutation.
$d = $word.ActiveDocument
